# Update "想去人数" (want-to-go count) values in column F across sheets
# 展览 (Exhibition), 演出 (Performance), 全部类型 (All types)

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 2034
$ws1.Range("F3").Value  = 652
$ws1.Range("F4").Value  = 1300
$ws1.Range("F6").Value  = 67
$ws1.Range("F8").Value  = 368
$ws1.Range("F9").Value  = 154
$ws1.Range("F11").Value = 931
$ws1.Range("F12").Value = 290
$ws1.Range("F13").Value = 154
$ws1.Range("F17").Value = 312
$ws1.Range("F18").Value = 730
$ws1.Range("F19").Value = 109
$ws1.Range("F20").Value = 689
$ws1.Range("F21").Value = 229
$ws1.Range("F22").Value = 59
$ws1.Range("F23").Value = 946
$ws1.Range("F24").Value = 402
$ws1.Range("F25").Value = 222
$ws1.Range("F26").Value = 71
$ws1.Range("F27").Value = 333
$ws1.Range("F30").Value = 442

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 273
$ws2.Range("F8").Value = 84

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 2034
$ws4.Range("F4").Value  = 652
$ws4.Range("F5").Value  = 1300
$ws4.Range("F8").Value  = 67
$ws4.Range("F10").Value = 368
$ws4.Range("F11").Value = 154
$ws4.Range("F13").Value = 931
$ws4.Range("F14").Value = 290
$ws4.Range("F15").Value = 154
$ws4.Range("F23").Value = 273
$ws4.Range("F24").Value = 313
$ws4.Range("F25").Value = 730
$ws4.Range("F26").Value = 109
$ws4.Range("F27").Value = 689
$ws4.Range("F28").Value = 229
$ws4.Range("F29").Value = 59
$ws4.Range("F30").Value = 946
$ws4.Range("F31").Value = 402
$ws4.Range("F32").Value = 84
$ws4.Range("F34").Value = 222
$ws4.Range("F35").Value = 71
$ws4.Range("F36").Value = 333
$ws4.Range("F42").Value = 442
